$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 9598.727999999999
$ws.Range("I70").Value = 20378
$ws.Range("J70").Value = 616
$ws.Range("K70").Value = 61134
$ws.Range("L70").Value = 1848
$ws.Range("M70").Value = -60864
$ws.Range("N70").Value = -2388

$ws.Range("H73").Value = 9598.727999999999
$ws.Range("I73").Value = 20378
$ws.Range("J73").Value = 616
$ws.Range("K73").Value = 61134
$ws.Range("L73").Value = 1848
$ws.Range("M73").Value = -60198
$ws.Range("N73").Value = -3720

$ws.Range("H138").Value = 3114.4517
$ws.Range("I138").Value = 1361.1666
$ws.Range("J138").Value = 3535.24
$ws.Range("K138").Value = 4083.4998
$ws.Range("L138").Value = 10605.72
$ws.Range("M138").Value = 1056.5002
$ws.Range("N138").Value = -20885.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18272.25
$ws.Range("I32").Value = 23002.646
$ws.Range("J32").Value = 4081.0625
$ws.Range("K32").Value = 23002.646
$ws.Range("L32").Value = 4081.0625
$ws.Range("M32").Value = -22715.646
$ws.Range("N32").Value = -4655.0625

$ws.Range("H63").Value = 2843154.8
$ws.Range("J63").Value = 15627003
$ws.Range("L63").Value = 15627003
$ws.Range("N63").Value = -15628375

$ws.Range("H66").Value = 2843154.8
$ws.Range("J66").Value = 15627003
$ws.Range("L66").Value = 78135015
$ws.Range("N66").Value = -78141879

$ws.Range("H80").Value = 41480
$ws.Range("J80").Value = 41480
$ws.Range("L80").Value = 41480
$ws.Range("N80").Value = -43476

$ws.Range("H83").Value = 41480
$ws.Range("J83").Value = 41480
$ws.Range("L83").Value = 124440
$ws.Range("N83").Value = -134424

$ws.Range("H122").Value = 2543.6072
$ws.Range("I122").Value = 2560.9
$ws.Range("J122").Value = 2500.375
$ws.Range("K122").Value = 7682.700000000001
$ws.Range("L122").Value = 7501.125
$ws.Range("M122").Value = -5232.700000000001
$ws.Range("N122").Value = -12401.125

$ws.Range("H132").Value = 104242.4
$ws.Range("I132").Value = 5303
$ws.Range("K132").Value = 15909
$ws.Range("M132").Value = -13379

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 8577
$ws.Range("J6").Value = 8577
$ws.Range("L6").Value = 8577
$ws.Range("N6").Value = -8803

$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620

$ws.Range("H107").Value = 582.2
$ws.Range("I107").Value = 602.75
$ws.Range("K107").Value = 602.75
$ws.Range("M107").Value = 1317.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22018.875
$ws.Range("I58").Value = 1357.6
$ws.Range("J58").Value = 56454.332
$ws.Range("K58").Value = 1357.6
$ws.Range("L58").Value = 56454.332
$ws.Range("M58").Value = -1154.6
$ws.Range("N58").Value = -56860.332

$ws.Range("H99").Value = 4434.7393
$ws.Range("I99").Value = 4137.375
$ws.Range("J99").Value = 4593.3335
$ws.Range("K99").Value = 4137.375
$ws.Range("L99").Value = 4593.3335
$ws.Range("M99").Value = -2639.375
$ws.Range("N99").Value = -7589.3335

$ws.Range("H126").Value = 4434.7393
$ws.Range("I126").Value = 4137.375
$ws.Range("J126").Value = 4593.3335
$ws.Range("K126").Value = 12412.125
$ws.Range("L126").Value = 13780.0005
$ws.Range("M126").Value = -9942.125
$ws.Range("N126").Value = -18720.0005

$ws.Range("H136").Value = 22018.875
$ws.Range("I136").Value = 1357.6
$ws.Range("J136").Value = 56454.332
$ws.Range("K136").Value = 4072.8
$ws.Range("L136").Value = 169362.996
$ws.Range("M136").Value = -1522.8
$ws.Range("N136").Value = -174462.996

$ws.Range("H139").Value = 9780
$ws.Range("J139").Value = 9780
$ws.Range("L139").Value = 9780
$ws.Range("N139").Value = -20060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3187.7234
$ws.Range("J68").Value = 3407.535
$ws.Range("L68").Value = 10222.605
$ws.Range("N68").Value = -11844.605

$ws.Range("H71").Value = 3187.7234
$ws.Range("J71").Value = 3407.535
$ws.Range("L71").Value = 30667.815
$ws.Range("N71").Value = -38779.815

$ws.Range("H107").Value = 4401.8623
$ws.Range("J107").Value = 1038.32
$ws.Range("L107").Value = 3114.96
$ws.Range("N107").Value = -6954.96

$ws.Range("H113").Value = 6433
$ws.Range("J113").Value = 751.75
$ws.Range("L113").Value = 2255.25
$ws.Range("N113").Value = -6595.25

$ws.Range("H131").Value = 137831.16
$ws.Range("I131").Value = 788.3333
$ws.Range("J131").Value = 150103.64
$ws.Range("K131").Value = 2364.9999
$ws.Range("L131").Value = 450310.92
$ws.Range("M131").Value = 2675.0001
$ws.Range("N131").Value = -460390.92

$ws.Range("H132").Value = 730.5333000000001
$ws.Range("I132").Value = 696.1818
$ws.Range("J132").Value = 825
$ws.Range("K132").Value = 6265.6362
$ws.Range("L132").Value = 7425
$ws.Range("M132").Value = -3735.6362
$ws.Range("N132").Value = -12485

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7167.8335
$ws.Range("I122").Value = 6503.5
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 19510.5
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -17060.5
$ws.Range("N122").Value = -27400

$ws.Range("H132").Value = 194111.12
$ws.Range("I132").Value = 260000.25
$ws.Range("J132").Value = 128222
$ws.Range("K132").Value = 780000.75
$ws.Range("L132").Value = 384666
$ws.Range("M132").Value = -777470.75
$ws.Range("N132").Value = -389726

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1394.2
$ws.Range("I22").Value = 1711.5714
$ws.Range("J22").Value = 653.6667
$ws.Range("K22").Value = 1711.5714
$ws.Range("L22").Value = 653.6667
$ws.Range("M22").Value = -1416.5714
$ws.Range("N22").Value = -1243.6667

$ws.Range("H27").Value = 1394.2
$ws.Range("I27").Value = 1711.5714
$ws.Range("J27").Value = 653.6667
$ws.Range("K27").Value = 1711.5714
$ws.Range("L27").Value = 653.6667
$ws.Range("M27").Value = -1604.5714
$ws.Range("N27").Value = -867.6667

$ws.Range("H46").Value = 885.7857
$ws.Range("I46").Value = 608.4167
$ws.Range("K46").Value = 608.4167
$ws.Range("M46").Value = -420.4167

$ws.Range("H61").Value = 4606.2085
$ws.Range("I61").Value = 2311.077
$ws.Range("J61").Value = 7318.636
$ws.Range("K61").Value = 2311.077
$ws.Range("L61").Value = 7318.636
$ws.Range("M61").Value = -2109.077
$ws.Range("N61").Value = -7722.636

$ws.Range("H113").Value = 4606.2085
$ws.Range("I113").Value = 2311.077
$ws.Range("J113").Value = 7318.636
$ws.Range("K113").Value = 2311.077
$ws.Range("L113").Value = 7318.636
$ws.Range("M113").Value = -141.0770000000002
$ws.Range("N113").Value = -11658.636

$ws.Range("H132").Value = 1750.4445
$ws.Range("I132").Value = 1280.76
$ws.Range("J132").Value = 2817.9092
$ws.Range("K132").Value = 3842.28
$ws.Range("L132").Value = 8453.7276
$ws.Range("M132").Value = -1312.28
$ws.Range("N132").Value = -13513.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1890.3793
$ws.Range("I132").Value = 1720.1428
$ws.Range("J132").Value = 2337.25
$ws.Range("K132").Value = 5160.428400000001
$ws.Range("L132").Value = 7011.75
$ws.Range("M132").Value = -2630.428400000001
$ws.Range("N132").Value = -12071.75
